# Update the build timestamp embedded in the version string from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.
#
# NOTE: reading a cell via .Value and feeding it straight back into another
# .Value assignment does not round-trip the text through this COM shim, so
# every cell's current text is read with .Text (which round-trips fine) and
# written back with .Value.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# About sheet: A2 (merged A2:G2) and A6 (merged A6:G6 region's anchor cell)
$a2 = $wsAbout.Range("A2").Text
$wsAbout.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6").Text
$wsAbout.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# Boundaries and methane sources sheet: S2:S11 hold the build_version string
for ($row = 2; $row -le 11; $row++) {
    $cell = $wsData.Cells.Item($row, 19)   # column S = 19
    $text = $cell.Text
    if ($text.Contains($oldStamp)) {
        $cell.Value = $text.Replace($oldStamp, $newStamp)
    }
}
